# Applies the edits described by the diff to cargaRGM_final.xlsx:
#  - A2: "126-2023-MPH/A" -> "41-2025-MPH/GM"
#  - B2: new cell value containing a truncated (200-char) copy of the long
#        description text (C2 keeps the full text, unchanged)
#  - D2: "03/01/2025" -> "03/09/2025"
#  - H2: "RESOLUCION 126-2023-MPH/A" -> "RESOLUCION 41-2025-MPH/GM"
#  - J2: "Archivo PDF (OCR:Sí) - 126-2023-MPH/A" -> "Documento 41-2025-MPH/GM"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full (unabridged) descriptive text that already lives in C2 and stays the same.
$fullText = "APROBAR la Liquidación Física y Financiera (con fines de verificación de cumplimiento de metas) del Proyecto: 2615348 “MEJORAMIENTO Y AMPLIACIÓN DE LOS SERVICIOS OPERATIVOS O MISIONALES INSTITUCIONALES EN LA SUBGERENCIA DE JUVENTUD, EDUCACIÓN Y DEPORTE DE LA MUNICIPALIDAD PROVINCIAL DE HUAMANGA DEL DISTRITO DE AYACUCHO DE LA PROVINCIA DE HUAMANGA DEL DEPARTAMENTO DE AYACUCHO”, ejecutado bajo la modalidad de Administración Directa, durante el ejercicio presupuestal del año 2023, de acuerdo al detalle que se consigna en la parte considerativa de la presente Resolución."

# Truncated text (first 200 characters) now stored in B2.
$truncatedText = $fullText.Substring(0, 200)

$ws.Range("A2").Value = "41-2025-MPH/GM"
$ws.Range("B2").Value = $truncatedText
$ws.Range("C2").Value = $fullText

# Keep the publication date as plain text (it was stored as text before the
# edit too), not as an auto-converted Excel date serial number. A leading
# apostrophe forces Excel to treat the value as text instead of parsing it
# into a date serial; the apostrophe itself is not part of the stored value.
$ws.Range("D2").Value = "'03/09/2025"
$ws.Range("D2").Style = "Normal"

$ws.Range("H2").Value = "RESOLUCION 41-2025-MPH/GM"
$ws.Range("J2").Value = "Documento 41-2025-MPH/GM"
